$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Assignment Title: Lab1 Understanding OOP Concepts with Code
# Examples" becomes "Assignment Title: Lab1 Part1 Understanding OOP Concepts
# with Code Examples" (a "Part1" designation is inserted after "Lab1").
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Assignment Title: Lab1 Understanding OOP Concepts with Code Examples",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Assignment Title: Lab1 Part1 Understanding OOP Concepts with Code Examples",
    2
)

# ---------------------------------------------------------------------------
# Change 2: the paragraph "Link to GitHub Lab1 repository:" (the one
# followed by the hyperlink to .../LAB1) becomes
# "Link to GitHub Lab1 Part1 repository:". Note there is a similarly worded
# paragraph "Link to GitHub Course repository:" earlier in the document that
# must stay untouched, so we match on the longer, unambiguous phrase.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Link to GitHub Lab1 repository:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Link to GitHub Lab1 Part1 repository:",
    2
)
